# Migrate "Biomass reactions" sheet to include a Compartment column in the
# biomass reaction table, and update the active-sheet/selection bookkeeping
# to match (Biomass reactions becomes the selected/active sheet instead of
# Database references).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Biomass reactions")

# Insert a new column C ("Compartment"), pushing the old C (Comments) and
# D (References) columns one to the right. Copying column B first means the
# newly inserted column inherits the same cell styles (header style / body
# style) that the rest of the table uses, rather than plain defaults.
$ws.Columns.Item(2).Copy()
$ws.Columns.Item(3).Insert()

# Match column B's display width for the new column.
$ws.Columns.Item(3).ColumnWidth = 22.330729166666668

# Fill in the new column's header + data cell.
$ws.Range("C1").Value = "Compartment"
$ws.Range("C2").Value = "c"

# Make "Biomass reactions" the active/selected sheet (was "Database
# references" before), with C3 selected.
$ws.Activate()
$ws.Range("C3").Select()
